# Add data for 2024-06-11
#
# A new day of incident data (2024-06-11) was appended to the source
# dataset. This pushes the running year-to-date "2024" total (column K)
# up on:
#   - the "Citywide Totals" sheet (one row per crime category + a Total row)
#   - the "By Neighborhood" sheet (one row per neighborhood + a Total row)
#   - every individual neighborhood sheet that had an incident that day
#     (one row per crime category + a Total row)
#
# All of the affected cells store plain numbers (no formulas), so this
# script simply overwrites column K for each affected row with the new
# running total.

$wb = $excel.ActiveWorkbook

# Each entry: worksheet index (tab order, 1-based), sheet name (for
# logging/readability only), row number, the new column-K value, and the
# previously-known value (used only as a sanity check before overwriting).
$changes = @(
    @{ Sheet = 1; Name = 'Citywide Totals'; Row = 2; New = 3376; Old = 3361 }
    @{ Sheet = 1; Name = 'Citywide Totals'; Row = 3; New = 3353; Old = 3331 }
    @{ Sheet = 1; Name = 'Citywide Totals'; Row = 4; New = 703; Old = 697 }
    @{ Sheet = 1; Name = 'Citywide Totals'; Row = 5; New = 221; Old = 220 }
    @{ Sheet = 1; Name = 'Citywide Totals'; Row = 6; New = 3943; Old = 3912 }
    @{ Sheet = 1; Name = 'Citywide Totals'; Row = 7; New = 11596; Old = 11521 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 7; New = 329; Old = 325 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 8; New = 770; Old = 759 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 9; New = 48; Old = 47 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 15; New = 118; Old = 117 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 19; New = 356; Old = 352 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 20; New = 266; Old = 264 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 21; New = 31; Old = 30 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 23; New = 112; Old = 111 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 25; New = 48; Old = 47 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 29; New = 606; Old = 604 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 33; New = 462; Old = 460 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 37; New = 403; Old = 402 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 41; New = 100; Old = 99 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 42; New = 415; Old = 412 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 43; New = 104; Old = 102 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 44; New = 108; Old = 106 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 46; New = 23; Old = 22 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 47; New = 63; Old = 62 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 51; New = 137; Old = 136 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 52; New = 316; Old = 314 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 53; New = 153; Old = 152 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 54; New = 229; Old = 225 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 58; New = 4; Old = 3 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 63; New = 40; Old = 38 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 65; New = 281; Old = 275 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 67; New = 452; Old = 449 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 72; New = 55; Old = 54 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 73; New = 103; Old = 102 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 78; New = 146; Old = 145 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 79; New = 299; Old = 298 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 83; New = 249; Old = 248 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 84; New = 82; Old = 81 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 85; New = 543; Old = 542 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 86; New = 77; Old = 76 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 88; New = 138; Old = 136 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 89; New = 154; Old = 153 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 90; New = 102; Old = 103 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 94; New = 143; Old = 142 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 95; New = 188; Old = 187 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 96; New = 142; Old = 141 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 97; New = 101; Old = 100 }
    @{ Sheet = 2; Name = 'By Neighborhood'; Row = 101; New = 11596; Old = 11521 }
    @{ Sheet = 4; Name = 'West Ridge'; Row = 3; New = 23; Old = 22 }
    @{ Sheet = 4; Name = 'West Ridge'; Row = 7; New = 142; Old = 141 }
    @{ Sheet = 5; Name = 'Auburn Gresham'; Row = 2; New = 120; Old = 117 }
    @{ Sheet = 5; Name = 'Auburn Gresham'; Row = 3; New = 99; Old = 98 }
    @{ Sheet = 5; Name = 'Auburn Gresham'; Row = 7; New = 329; Old = 325 }
    @{ Sheet = 7; Name = 'Uptown'; Row = 3; New = 50; Old = 49 }
    @{ Sheet = 7; Name = 'Uptown'; Row = 7; New = 154; Old = 153 }
    @{ Sheet = 8; Name = 'South Shore'; Row = 6; New = 126; Old = 125 }
    @{ Sheet = 8; Name = 'South Shore'; Row = 7; New = 543; Old = 542 }
    @{ Sheet = 9; Name = 'Little Village'; Row = 6; New = 127; Old = 125 }
    @{ Sheet = 9; Name = 'Little Village'; Row = 7; New = 316; Old = 314 }
    @{ Sheet = 11; Name = 'Logan Square'; Row = 6; New = 78; Old = 77 }
    @{ Sheet = 11; Name = 'Logan Square'; Row = 7; New = 153; Old = 152 }
    @{ Sheet = 12; Name = 'Austin'; Row = 2; New = 222; Old = 221 }
    @{ Sheet = 12; Name = 'Austin'; Row = 3; New = 232; Old = 230 }
    @{ Sheet = 12; Name = 'Austin'; Row = 4; New = 43; Old = 42 }
    @{ Sheet = 12; Name = 'Austin'; Row = 6; New = 252; Old = 245 }
    @{ Sheet = 12; Name = 'Austin'; Row = 7; New = 770; Old = 759 }
    @{ Sheet = 13; Name = 'South Chicago'; Row = 6; New = 60; Old = 59 }
    @{ Sheet = 13; Name = 'South Chicago'; Row = 7; New = 249; Old = 248 }
    @{ Sheet = 14; Name = 'Garfield Park'; Row = 3; New = 173; Old = 171 }
    @{ Sheet = 14; Name = 'Garfield Park'; Row = 7; New = 462; Old = 460 }
    @{ Sheet = 15; Name = 'West Pullman'; Row = 3; New = 67; Old = 66 }
    @{ Sheet = 15; Name = 'West Pullman'; Row = 7; New = 188; Old = 187 }
    @{ Sheet = 16; Name = 'Grand Crossing'; Row = 2; New = 108; Old = 109 }
    @{ Sheet = 16; Name = 'Grand Crossing'; Row = 6; New = 124; Old = 122 }
    @{ Sheet = 16; Name = 'Grand Crossing'; Row = 7; New = 403; Old = 402 }
    @{ Sheet = 17; Name = 'New City'; Row = 3; New = 73; Old = 71 }
    @{ Sheet = 17; Name = 'New City'; Row = 4; New = 9; Old = 8 }
    @{ Sheet = 17; Name = 'New City'; Row = 6; New = 109; Old = 106 }
    @{ Sheet = 17; Name = 'New City'; Row = 7; New = 281; Old = 275 }
    @{ Sheet = 21; Name = 'North Lawndale'; Row = 2; New = 140; Old = 139 }
    @{ Sheet = 21; Name = 'North Lawndale'; Row = 3; New = 147; Old = 146 }
    @{ Sheet = 21; Name = 'North Lawndale'; Row = 6; New = 132; Old = 131 }
    @{ Sheet = 21; Name = 'North Lawndale'; Row = 7; New = 452; Old = 449 }
    @{ Sheet = 22; Name = 'South Deering'; Row = 3; New = 31; Old = 30 }
    @{ Sheet = 22; Name = 'South Deering'; Row = 7; New = 82; Old = 81 }
    @{ Sheet = 24; Name = 'Loop'; Row = 2; New = 41; Old = 40 }
    @{ Sheet = 24; Name = 'Loop'; Row = 3; New = 68; Old = 67 }
    @{ Sheet = 24; Name = 'Loop'; Row = 6; New = 107; Old = 105 }
    @{ Sheet = 24; Name = 'Loop'; Row = 7; New = 229; Old = 225 }
    @{ Sheet = 25; Name = 'Englewood'; Row = 3; New = 206; Old = 205 }
    @{ Sheet = 25; Name = 'Englewood'; Row = 6; New = 185; Old = 184 }
    @{ Sheet = 25; Name = 'Englewood'; Row = 7; New = 606; Old = 604 }
    @{ Sheet = 27; Name = 'Chatham'; Row = 2; New = 121; Old = 119 }
    @{ Sheet = 27; Name = 'Chatham'; Row = 3; New = 95; Old = 93 }
    @{ Sheet = 27; Name = 'Chatham'; Row = 7; New = 356; Old = 352 }
    @{ Sheet = 28; Name = 'Irving Park'; Row = 2; New = 20; Old = 19 }
    @{ Sheet = 28; Name = 'Irving Park'; Row = 6; New = 51; Old = 50 }
    @{ Sheet = 28; Name = 'Irving Park'; Row = 7; New = 108; Old = 106 }
    @{ Sheet = 31; Name = 'Hermosa'; Row = 6; New = 42; Old = 41 }
    @{ Sheet = 31; Name = 'Hermosa'; Row = 7; New = 100; Old = 99 }
    @{ Sheet = 32; Name = 'Humboldt Park'; Row = 2; New = 110; Old = 109 }
    @{ Sheet = 32; Name = 'Humboldt Park'; Row = 3; New = 133; Old = 131 }
    @{ Sheet = 32; Name = 'Humboldt Park'; Row = 7; New = 415; Old = 412 }
    @{ Sheet = 35; Name = 'Rogers Park'; Row = 2; New = 44; Old = 43 }
    @{ Sheet = 35; Name = 'Rogers Park'; Row = 7; New = 146; Old = 145 }
    @{ Sheet = 38; Name = 'Jefferson Park'; Row = 2; New = 9; Old = 8 }
    @{ Sheet = 38; Name = 'Jefferson Park'; Row = 7; New = 23; Old = 22 }
    @{ Sheet = 39; Name = 'Douglas'; Row = 6; New = 29; Old = 28 }
    @{ Sheet = 39; Name = 'Douglas'; Row = 7; New = 112; Old = 111 }
    @{ Sheet = 41; Name = 'Chinatown'; Row = 6; New = 17; Old = 16 }
    @{ Sheet = 41; Name = 'Chinatown'; Row = 7; New = 31; Old = 30 }
    @{ Sheet = 42; Name = 'Roseland'; Row = 3; New = 103; Old = 102 }
    @{ Sheet = 42; Name = 'Roseland'; Row = 7; New = 299; Old = 298 }
    @{ Sheet = 43; Name = 'Near South Side'; Row = 3; New = 24; Old = 23 }
    @{ Sheet = 43; Name = 'Near South Side'; Row = 6; New = 27; Old = 28 }
    @{ Sheet = 44; Name = 'Chicago Lawn'; Row = 2; New = 93; Old = 92 }
    @{ Sheet = 44; Name = 'Chicago Lawn'; Row = 3; New = 75; Old = 74 }
    @{ Sheet = 44; Name = 'Chicago Lawn'; Row = 7; New = 266; Old = 264 }
    @{ Sheet = 51; Name = 'West Loop'; Row = 2; New = 40; Old = 39 }
    @{ Sheet = 51; Name = 'West Loop'; Row = 7; New = 143; Old = 142 }
    @{ Sheet = 52; Name = 'East Side'; Row = 6; New = 8; Old = 7 }
    @{ Sheet = 52; Name = 'East Side'; Row = 7; New = 48; Old = 47 }
    @{ Sheet = 53; Name = 'Kenwood'; Row = 2; New = 19; Old = 18 }
    @{ Sheet = 53; Name = 'Kenwood'; Row = 7; New = 63; Old = 62 }
    @{ Sheet = 54; Name = 'Brighton Park'; Row = 2; New = 42; Old = 41 }
    @{ Sheet = 54; Name = 'Brighton Park'; Row = 7; New = 118; Old = 117 }
    @{ Sheet = 61; Name = 'Avalon Park'; Row = 6; New = 11; Old = 10 }
    @{ Sheet = 61; Name = 'Avalon Park'; Row = 7; New = 48; Old = 47 }
    @{ Sheet = 62; Name = 'Portage Park'; Row = 6; New = 42; Old = 41 }
    @{ Sheet = 62; Name = 'Portage Park'; Row = 7; New = 103; Old = 102 }
    @{ Sheet = 65; Name = 'West Town'; Row = 6; New = 62; Old = 61 }
    @{ Sheet = 65; Name = 'West Town'; Row = 7; New = 101; Old = 100 }
    @{ Sheet = 68; Name = 'United Center'; Row = 3; New = 37; Old = 35 }
    @{ Sheet = 68; Name = 'United Center'; Row = 7; New = 138; Old = 136 }
    @{ Sheet = 72; Name = 'Streeterville'; Row = 4; New = 29; Old = 28 }
    @{ Sheet = 72; Name = 'Streeterville'; Row = 7; New = 77; Old = 76 }
    @{ Sheet = 74; Name = 'Washington Heights'; Row = 4; New = 9; Old = 10 }
    @{ Sheet = 74; Name = 'Washington Heights'; Row = 7; New = 102; Old = 103 }
    @{ Sheet = 75; Name = 'Little Italy, UIC'; Row = 5; New = 3; Old = 2 }
    @{ Sheet = 75; Name = 'Little Italy, UIC'; Row = 7; New = 137; Old = 136 }
    @{ Sheet = 79; Name = 'Hyde Park'; Row = 2; New = 21; Old = 20 }
    @{ Sheet = 79; Name = 'Hyde Park'; Row = 6; New = 44; Old = 43 }
    @{ Sheet = 79; Name = 'Hyde Park'; Row = 7; New = 104; Old = 102 }
    @{ Sheet = 82; Name = 'Old Town'; Row = 4; New = 4; Old = 3 }
    @{ Sheet = 82; Name = 'Old Town'; Row = 7; New = 55; Old = 54 }
    @{ Sheet = 97; Name = 'Millenium Park'; Row = 6; New = 4; Old = 3 }
    @{ Sheet = 97; Name = 'Millenium Park'; Row = 7; New = 4; Old = 3 }
)

$mismatchCount = 0
$updateCount = 0
$sheetsTouched = @{}

foreach ($item in $changes) {
    $ws = $wb.Worksheets.Item($item.Sheet)
    $cell = $ws.Cells.Item($item.Row, 11)
    $current = $cell.Value2

    if ($current -ne $item.Old) {
        $mismatchCount++
        Write-Host "WARNING:" $item.Name "row" $item.Row "expected old value" $item.Old "but found" $current
    }

    $cell.Value = $item.New
    $updateCount++
    $sheetsTouched[$item.Sheet] = $true
}

Write-Host "Updated" $updateCount "cells across" $sheetsTouched.Keys.Count "sheets."
if ($mismatchCount -gt 0) {
    Write-Host "Mismatches encountered:" $mismatchCount
} else {
    Write-Host "All prior values matched expectations."
}
